$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (adegraphics) was missing a tutorial/vignette entry - mark it "NA".
$ws.Range("E5").Value = "NA"

# New row 12: add the pcaMethods package to the comparison table.
$ws.Range("A12").Value = "pcaMethods"
$ws.Range("D12").Value = "ggplot"
$ws.Range("E12").Value = "The pcaMethods Package"

# Hyperlink the package name to its package page, and the write-up to its vignette,
# then restore the "Hyperlink" cell style (Hyperlinks.Add mints its own style xf).
$ws.Hyperlinks.Add($ws.Range("A12"), "https://www.bioconductor.org/packages/release/bioc/html/pcaMethods.html", "", "", "pcaMethods")
$ws.Range("A12").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E12"), "https://www.bioconductor.org/packages/release/bioc/vignettes/pcaMethods/inst/doc/pcaMethods.pdf", "", "", "The pcaMethods Package")
$ws.Range("E12").Style = "Hyperlink"

# Move the active selection to match the author's final cursor position.
$ws.Range("E5").Select()
